# Natmi following Dr Hou advice
# Update LR-pair statistics for Sema4d-Plxnb2 (recomputed ligand/receptor cell counts and derived metrics)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 55.59510733333332
$ws.Cells.Item(2, 8).Value = 166.785322
$ws.Cells.Item(2, 9).Value = 0.4537221086682116
$ws.Cells.Item(2, 10).Value = 0.4537221086682116
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 10.39091433333333
$ws.Cells.Item(2, 14).Value = 31.172743
$ws.Cells.Item(2, 15).Value = 0.0835098648954196
$ws.Cells.Item(2, 16).Value = 0.0835098648954196
$ws.Cells.Item(2, 17).Value = 577.6839976531384
$ws.Cells.Item(2, 18).Value = 5199.155978878245
$ws.Cells.Item(2, 19).Value = 0.03789027199494725
$ws.Cells.Item(2, 20).Value = 0.03789027199494724

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 55.59510733333332
$ws.Cells.Item(3, 8).Value = 166.785322
$ws.Cells.Item(3, 9).Value = 0.4537221086682116
$ws.Cells.Item(3, 10).Value = 0.4537221086682116
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 35.71561933333334
$ws.Cells.Item(3, 14).Value = 107.146858
$ws.Cells.Item(3, 15).Value = 0.287039855156433
$ws.Cells.Item(3, 16).Value = 0.287039855156433
$ws.Cells.Item(3, 17).Value = 1985.613690313142
$ws.Cells.Item(3, 18).Value = 17870.52321281828
$ws.Cells.Item(3, 19).Value = 0.1302363283533948
$ws.Cells.Item(3, 20).Value = 0.1302363283533948

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 55.59510733333332
$ws.Cells.Item(4, 8).Value = 166.785322
$ws.Cells.Item(4, 9).Value = 0.4537221086682116
$ws.Cells.Item(4, 10).Value = 0.4537221086682116
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 51.87044599999999
$ws.Cells.Item(4, 14).Value = 155.611338
$ws.Cells.Item(4, 15).Value = 0.4168732219867682
$ws.Cells.Item(4, 16).Value = 0.4168732219867682
$ws.Cells.Item(4, 17).Value = 2883.74301279787
$ws.Cells.Item(4, 18).Value = 25953.68711518083
$ws.Cells.Item(4, 19).Value = 0.1891445973271479
$ws.Cells.Item(4, 20).Value = 0.1891445973271479

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 55.59510733333332
$ws.Cells.Item(5, 8).Value = 166.785322
$ws.Cells.Item(5, 9).Value = 0.4537221086682116
$ws.Cells.Item(5, 10).Value = 0.4537221086682116
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 26.450408
$ws.Cells.Item(5, 14).Value = 79.351224
$ws.Cells.Item(5, 15).Value = 0.2125770579613792
$ws.Cells.Item(5, 16).Value = 0.2125770579613792
$ws.Cells.Item(5, 17).Value = 1470.513271770458
$ws.Cells.Item(5, 18).Value = 13234.61944593413
$ws.Cells.Item(5, 19).Value = 0.09645091099272161
$ws.Cells.Item(5, 20).Value = 0.09645091099272161

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.654659
$ws.Cells.Item(6, 8).Value = 4.963977
$ws.Cells.Item(6, 9).Value = 0.01350398275347337
$ws.Cells.Item(6, 10).Value = 0.01350398275347337
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 10.39091433333333
$ws.Cells.Item(6, 14).Value = 31.172743
$ws.Cells.Item(6, 15).Value = 0.0835098648954196
$ws.Cells.Item(6, 16).Value = 0.0835098648954196
$ws.Cells.Item(6, 17).Value = 17.193419919879
$ws.Cells.Item(6, 18).Value = 154.740779278911
$ws.Cells.Item(6, 19).Value = 0.001127715775292638
$ws.Cells.Item(6, 20).Value = 0.001127715775292638

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.654659
$ws.Cells.Item(7, 8).Value = 4.963977
$ws.Cells.Item(7, 9).Value = 0.01350398275347337
$ws.Cells.Item(7, 10).Value = 0.01350398275347337
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 35.71561933333334
$ws.Cells.Item(7, 14).Value = 107.146858
$ws.Cells.Item(7, 15).Value = 0.287039855156433
$ws.Cells.Item(7, 16).Value = 0.287039855156433
$ws.Cells.Item(7, 17).Value = 59.097170970474
$ws.Cells.Item(7, 18).Value = 531.874538734266
$ws.Cells.Item(7, 19).Value = 0.003876181253591966
$ws.Cells.Item(7, 20).Value = 0.003876181253591966

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.654659
$ws.Cells.Item(8, 8).Value = 4.963977
$ws.Cells.Item(8, 9).Value = 0.01350398275347337
$ws.Cells.Item(8, 10).Value = 0.01350398275347337
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 51.87044599999999
$ws.Cells.Item(8, 14).Value = 155.611338
$ws.Cells.Item(8, 15).Value = 0.4168732219867682
$ws.Cells.Item(8, 16).Value = 0.4168732219867682
$ws.Cells.Item(8, 17).Value = 85.82790030791398
$ws.Cells.Item(8, 18).Value = 772.4511027712259
$ws.Cells.Item(8, 19).Value = 0.005629448800094194
$ws.Cells.Item(8, 20).Value = 0.005629448800094194

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.654659
$ws.Cells.Item(9, 8).Value = 4.963977
$ws.Cells.Item(9, 9).Value = 0.01350398275347337
$ws.Cells.Item(9, 10).Value = 0.01350398275347337
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 26.450408
$ws.Cells.Item(9, 14).Value = 79.351224
$ws.Cells.Item(9, 15).Value = 0.2125770579613792
$ws.Cells.Item(9, 16).Value = 0.2125770579613792
$ws.Cells.Item(9, 17).Value = 43.766405650872
$ws.Cells.Item(9, 18).Value = 393.897650857848
$ws.Cells.Item(9, 19).Value = 0.002870636924494574
$ws.Cells.Item(9, 20).Value = 0.002870636924494574

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 62.05924166666666
$ws.Cells.Item(10, 8).Value = 186.177725
$ws.Cells.Item(10, 9).Value = 0.5064771225734745
$ws.Cells.Item(10, 10).Value = 0.5064771225734744
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 10.39091433333333
$ws.Cells.Item(10, 14).Value = 31.172743
$ws.Cells.Item(10, 15).Value = 0.0835098648954196
$ws.Cells.Item(10, 16).Value = 0.0835098648954196
$ws.Cells.Item(10, 17).Value = 644.8522637499638
$ws.Cells.Item(10, 18).Value = 5803.670373749674
$ws.Cells.Item(10, 19).Value = 0.04229583607873173
$ws.Cells.Item(10, 20).Value = 0.04229583607873173

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 62.05924166666666
$ws.Cells.Item(11, 8).Value = 186.177725
$ws.Cells.Item(11, 9).Value = 0.5064771225734745
$ws.Cells.Item(11, 10).Value = 0.5064771225734744
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 35.71561933333334
$ws.Cells.Item(11, 14).Value = 107.146858
$ws.Cells.Item(11, 15).Value = 0.287039855156433
$ws.Cells.Item(11, 16).Value = 0.287039855156433
$ws.Cells.Item(11, 17).Value = 2216.484251482005
$ws.Cells.Item(11, 18).Value = 19948.35826333805
$ws.Cells.Item(11, 19).Value = 0.1453791199035371
$ws.Cells.Item(11, 20).Value = 0.1453791199035371

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 62.05924166666666
$ws.Cells.Item(12, 8).Value = 186.177725
$ws.Cells.Item(12, 9).Value = 0.5064771225734745
$ws.Cells.Item(12, 10).Value = 0.5064771225734744
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 51.87044599999999
$ws.Cells.Item(12, 14).Value = 155.611338
$ws.Cells.Item(12, 15).Value = 0.4168732219867682
$ws.Cells.Item(12, 16).Value = 0.4168732219867682
$ws.Cells.Item(12, 17).Value = 3219.040543671782
$ws.Cells.Item(12, 18).Value = 28971.36489304604
$ws.Cells.Item(12, 19).Value = 0.2111367499497916
$ws.Cells.Item(12, 20).Value = 0.2111367499497916

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 62.05924166666666
$ws.Cells.Item(13, 8).Value = 186.177725
$ws.Cells.Item(13, 9).Value = 0.5064771225734745
$ws.Cells.Item(13, 10).Value = 0.5064771225734744
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 26.450408
$ws.Cells.Item(13, 14).Value = 79.351224
$ws.Cells.Item(13, 15).Value = 0.2125770579613792
$ws.Cells.Item(13, 16).Value = 0.2125770579613792
$ws.Cells.Item(13, 17).Value = 1641.492262253933
$ws.Cells.Item(13, 18).Value = 14773.4303602854
$ws.Cells.Item(13, 19).Value = 0.1076654166414141
$ws.Cells.Item(13, 20).Value = 0.107665416641414

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 3.222176333333334
$ws.Cells.Item(14, 8).Value = 9.666529
$ws.Cells.Item(14, 9).Value = 0.02629678600484052
$ws.Cells.Item(14, 10).Value = 0.02629678600484052
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 10.39091433333333
$ws.Cells.Item(14, 14).Value = 31.172743
$ws.Cells.Item(14, 15).Value = 0.0835098648954196
$ws.Cells.Item(14, 16).Value = 0.0835098648954196
$ws.Cells.Item(14, 17).Value = 33.48135824656078
$ws.Cells.Item(14, 18).Value = 301.332224219047
$ws.Cells.Item(14, 19).Value = 0.002196041046447993
$ws.Cells.Item(14, 20).Value = 0.002196041046447992

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 3.222176333333334
$ws.Cells.Item(15, 8).Value = 9.666529
$ws.Cells.Item(15, 9).Value = 0.02629678600484052
$ws.Cells.Item(15, 10).Value = 0.02629678600484052
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 35.71561933333334
$ws.Cells.Item(15, 14).Value = 107.146858
$ws.Cells.Item(15, 15).Value = 0.287039855156433
$ws.Cells.Item(15, 16).Value = 0.287039855156433
$ws.Cells.Item(15, 17).Value = 115.0820233462091
$ws.Cells.Item(15, 18).Value = 1035.738210115882
$ws.Cells.Item(15, 19).Value = 0.007548225645909136
$ws.Cells.Item(15, 20).Value = 0.007548225645909136

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 3.222176333333334
$ws.Cells.Item(16, 8).Value = 9.666529
$ws.Cells.Item(16, 9).Value = 0.02629678600484052
$ws.Cells.Item(16, 10).Value = 0.02629678600484052
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 51.87044599999999
$ws.Cells.Item(16, 14).Value = 155.611338
$ws.Cells.Item(16, 15).Value = 0.4168732219867682
$ws.Cells.Item(16, 16).Value = 0.4168732219867682
$ws.Cells.Item(16, 17).Value = 167.1357235006446
$ws.Cells.Item(16, 18).Value = 1504.221511505802
$ws.Cells.Item(16, 19).Value = 0.01096242590973442
$ws.Cells.Item(16, 20).Value = 0.01096242590973442

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 3.222176333333334
$ws.Cells.Item(17, 8).Value = 9.666529
$ws.Cells.Item(17, 9).Value = 0.02629678600484052
$ws.Cells.Item(17, 10).Value = 0.02629678600484052
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 26.450408
$ws.Cells.Item(17, 14).Value = 79.351224
$ws.Cells.Item(17, 15).Value = 0.2125770579613792
$ws.Cells.Item(17, 16).Value = 0.2125770579613792
$ws.Cells.Item(17, 17).Value = 85.22787866461067
$ws.Cells.Item(17, 18).Value = 767.050907981496
$ws.Cells.Item(17, 19).Value = 0.005590093402748969
$ws.Cells.Item(17, 20).Value = 0.005590093402748969

